$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '43.953.53'
$ws.Range('E2').Value = '  -1.13%  '

# Row 3
$ws.Range('D3').Value = '2.239.23'
$ws.Range('E3').Value = '  -1.94%  '

# Row 4
$ws.Range('E4').Value = '  +0.07%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.64'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.92%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '98.91'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -7.13%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.573'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -3.37%  '

# Row 8
$ws.Range('E8').Value = '  +0.03%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.530'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -7.64%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.04'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -7.42%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0821'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.00%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.33'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -7.67%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.104'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.89%  '

# Row 14
$ws.Range('D14').Value = '2.579.31'
$ws.Range('E14').Value = '  -2.06%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.839'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -5.71%  '

# Row 16
$ws.Range('D16').Value = '2.237.79'
$ws.Range('E16').Value = '  -1.91%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.90'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -5.39%  '

# Row 18
$ws.Range('D18').Value = '43.821.24'
$ws.Range('E18').Value = '  -1.27%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.08'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -8.06%  '

# Row 20
$ws.Range('D20').Value = '0.0₃0976'
$ws.Range('E20').Value = '  -3.12%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.31'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.89%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '65.46'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.78%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '235.01'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.98%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.97'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -7.89%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.01'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -8.98%  '

# Row 26
$ws.Range('E26').Value = '  +0.26%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.13'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.18%  '

# Row 28
$ws.Range('E28').Value = '  -4.65%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '36.25'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -6.08%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.95'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -9.60%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '20.00'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.48%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '154.92'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.08%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0832'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -6.83%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.29'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.64%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.66'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.93%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.89'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -8.69%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.107'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -8.67%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.117'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.65%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '15.45'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.15%  '

# Row 40
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.99'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -11.24%  '

# Row 41
$ws.Range('B41').Value = 'NEARProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.50'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -11.92%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0306'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -7.24%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.01'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.21%  '

# Row 44
$ws.Range('D44').Value = '1.697.17'
$ws.Range('E44').Value = '  -4.71%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '82.09'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -6.00%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.194'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -7.12%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '5.16'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -6.64%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '101.40'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.25%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '71.21'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -5.20%  '

# Row 50
$ws.Range('B50').Value = 'MultiversX'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '56.04'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -7.64%  '

# Row 51
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.59'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -7.82%  '
